$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "{edit:linked}"/"{edit:unlinked}" template cells (text unchanged,
# but kept here to be explicit / ensure correct shared-string state)
$ws.Range("P2").Value = "{edit:linked}"
$ws.Range("P3").Value = "{edit:unlinked}"

# Column header template cells get re-keyed from positional child[0..2]
# references to named child[FIELD] references
$ws.Range("A7").Value = "{child[ID]:linked}"
$ws.Range("B7").Value = "{child[DESCRIPTION]:doc_link}"
$ws.Range("C7").Value = "{child[DATE]}"

# Update the active selection on the sheet to match the new authoring target
$ws.Range("C7").Select()
